# "Actualización 11 de Mayo - Mañana"
# Updates the "Rescatables" sheet with the refreshed list of students that
# still need to take the make-up ("rescate") exam. The list grew from 7 to
# 16 students and the first group (2ALCV, "Reprobadas"=2) that used to be at
# the top got re-sorted after new 2ALCV/2AEV groups, and a new trailing block
# of "Reprobadas"=1 rows was appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$materia = "LECTURA, EXPRESIÓN ORAL Y ESCRITA II"

# Columns: A=Mat, B=Paterno, C=Materno, D=Nombres, E=Nombre_Largo, F=Grupo, G=Reprobadas
$data = @(
  @(20330051920359, "CASTILLO",  "ROMERO",    "KARLA JOVANA",      "2ALCV", 2),
  @(20330051920091, "OLMOS",     "CASTRO",    "ANGEL",             "2ALCV", 2),
  @(20330051920191, "ARIAS",     "BARRAGAN",  "ANALI",             "2ALCV", 2),
  @(20330051920202, "DE JESUS",  "CASTILLO",  "ITZEL",             "2ALCV", 2),
  @(20330051920218, "DE LA TEJA","RAMIREZ",   "ALISSON FERNANDA",  "2ALCV", 2),
  @(20330051920274, "MARTINEZ",  "RODRIGUEZ", "DANIEL ELEAZAR",    "2APV",  2),
  @(20330051920278, "RAMOS",     "XOTLANIHUA","MARCO JOSAFAT",     "2APV",  2),
  @(20330051920282, "TINOCO",    "RAMOS",     "ERNESTO",           "2APV",  2),
  @(20330051920041, "AVENDAÑO",  "SANCHEZ",   "AXEL JESUS",        "2AEV",  1),
  @(20330051920075, "CONTRERAS", "GARCIA",    "JORGE HUMBERTO",    "2AEV",  1),
  @(20330051920062, "ROMERO",    "REYES",     "AMANDA MICHEL",     "2AEV",  1),
  @(20330051920201, "GOMEZ",     "GONZALEZ",  "YAMILET",           "2ALCV", 1),
  @(20330051920210, "ROBLES",    "CASTILLO",  "JULIO CESAR",       "2ALCV", 1),
  @(20330051920265, "CORONA",    "HERNANDEZ", "GUADALUPE",         "2APV",  1),
  @(20330051920268, "GONZALEZ",  "FLORES",    "JESUS HUMBERTO",    "2APV",  1),
  @(20330051920281, "TELLEZ",    "OFICIAL",   "MARISOL",           "2APV",  1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $rec = $data[$i]

  $ws.Cells.Item($row, 1).Value = $rec[0]
  $ws.Cells.Item($row, 2).Value = $rec[1]
  $ws.Cells.Item($row, 3).Value = $rec[2]
  $ws.Cells.Item($row, 4).Value = $rec[3]
  $ws.Cells.Item($row, 5).Value = $materia
  $ws.Cells.Item($row, 6).Value = $rec[4]
  $ws.Cells.Item($row, 7).Value = $rec[5]
}
